# feat: Ui otimization + login page rework
# Add a new "VehicleGroup" column (K) to the agendamentos template and
# move the active selection to the new empty cell below the header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell K1 - mirror the look & feel of the existing header cells
# (bold, centered, wrapped text) by copying J1's formatting, then set the text.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("K1").Value = "VehicleGroup"

# Give the new column a sensible custom width.
$ws.Columns.Item(11).ColumnWidth = 14.333333333333334

# Header row grows a bit taller to match the rest of the table.
$ws.Rows.Item(1).RowHeight = 29

# Move the active selection onto the new column, ready for data entry.
$ws.Range("K2").Select()
